$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "student-id"
$ws.Range("C1").Value = "counsl Spring"
$ws.Range("D1").Value = "counsl Fall"
$ws.Range("B1").Value = "honors"
$ws.Range("E1").Value = "advising"
$ws.Range("F1").Value = "major"

$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 1

$headerRange = $ws.Range("A1:F1")
$headerRange.WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

$ws.Range("F2").Select() | Out-Null
